# Append: 2026-02-09 07:10 JST
# The scraper re-ran and refreshed the "acquired at" timestamp (column A,
# "取得日時") for every row currently on the listing sheet ("ランサーズ"),
# from 2026-02-09 07:01:52 to 2026-02-09 07:10:29. No other data changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2026-02-09 07:01:52"
$newTimestamp = "2026-02-09 07:10:29"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Text -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
